# Applies the scheduled-runner Chocobo Leve Profit price/profit refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H128").Value2 = 41846
$ws.Range("J128").Value2 = 41846
$ws.Range("L128").Value2 = 41846
$ws.Range("N128").Value2 = -51806

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value2 = 1000
$ws.Range("I2").Value2 = 1000
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 1000
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = -887
$ws.Range("N2").Value2 = $null

$ws.Range("H10").Value2 = 10555.6
$ws.Range("I10").Value2 = 2500
$ws.Range("J10").Value2 = 15926
$ws.Range("K10").Value2 = 2500
$ws.Range("L10").Value2 = 15926
$ws.Range("M10").Value2 = -2330
$ws.Range("N10").Value2 = -16266

$ws.Range("H25").Value2 = 5171.1665
$ws.Range("I25").Value2 = 1757
$ws.Range("J25").Value2 = 11999.5
$ws.Range("K25").Value2 = 1757
$ws.Range("L25").Value2 = 11999.5
$ws.Range("M25").Value2 = -1355
$ws.Range("N25").Value2 = -12803.5

$ws.Range("H27").Value2 = 30000
$ws.Range("J27").Value2 = 30000
$ws.Range("L27").Value2 = 30000
$ws.Range("N27").Value2 = -30368

$ws.Range("H30").Value2 = 7899.6
$ws.Range("I30").Value2 = 6666.3335
$ws.Range("J30").Value2 = 9749.5
$ws.Range("K30").Value2 = 6666.3335
$ws.Range("L30").Value2 = 9749.5
$ws.Range("M30").Value2 = -6516.3335
$ws.Range("N30").Value2 = -10049.5

$ws.Range("H116").Value2 = 1000
$ws.Range("I116").Value2 = 1000
$ws.Range("J116").Value2 = 0
$ws.Range("K116").Value2 = 1000
$ws.Range("L116").Value2 = 0
$ws.Range("M116").Value2 = 1294
$ws.Range("N116").Value2 = $null

$ws.Range("H122").Value2 = 3892.4
$ws.Range("I122").Value2 = 1115.5
$ws.Range("K122").Value2 = 3346.5
$ws.Range("M122").Value2 = -896.5

$ws.Range("H138").Value2 = 77490
$ws.Range("J138").Value2 = 77490
$ws.Range("L138").Value2 = 77490
$ws.Range("N138").Value2 = -87770

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value2 = 1000
$ws.Range("I3").Value2 = 1000
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 1000
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = -886
$ws.Range("N3").Value2 = $null

$ws.Range("H75").Value2 = 4000
$ws.Range("I75").Value2 = 4000
$ws.Range("J75").Value2 = 0
$ws.Range("K75").Value2 = 4000
$ws.Range("L75").Value2 = 0
$ws.Range("M75").Value2 = -3064
$ws.Range("N75").Value2 = $null

$ws.Range("H78").Value2 = 4000
$ws.Range("I78").Value2 = 4000
$ws.Range("J78").Value2 = 0
$ws.Range("K78").Value2 = 12000
$ws.Range("L78").Value2 = 0
$ws.Range("M78").Value2 = -7320
$ws.Range("N78").Value2 = $null

$ws = $wb.Worksheets("CRP")
$ws.Range("H14").Value2 = 10255.5
$ws.Range("J14").Value2 = 10255.5
$ws.Range("L14").Value2 = 10255.5
$ws.Range("N14").Value2 = -10595.5

$ws.Range("H110").Value2 = 41980
$ws.Range("J110").Value2 = 41980
$ws.Range("L110").Value2 = 41980
$ws.Range("N110").Value2 = -50160

$ws.Range("H111").Value2 = 40000
$ws.Range("J111").Value2 = 40000
$ws.Range("L111").Value2 = 40000
$ws.Range("N111").Value2 = -48180

$ws.Range("H112").Value2 = 30082.295
$ws.Range("J112").Value2 = 30082.295
$ws.Range("L112").Value2 = 30082.295
$ws.Range("N112").Value2 = -33036.295

$ws.Range("H122").Value2 = 15000
$ws.Range("I122").Value2 = 0
$ws.Range("K122").Value2 = 0
$ws.Range("M122").Value2 = $null

$ws = $wb.Worksheets("CUL")
$ws.Range("H47").Value2 = 2049.25
$ws.Range("I47").Value2 = 399
$ws.Range("K47").Value2 = 1197
$ws.Range("M47").Value2 = -766

$ws.Range("H69").Value2 = 5840
$ws.Range("I69").Value2 = 866.6667
$ws.Range("J69").Value2 = 8326.666999999999
$ws.Range("K69").Value2 = 2600.0001
$ws.Range("L69").Value2 = 24980.001
$ws.Range("M69").Value2 = -1789.0001
$ws.Range("N69").Value2 = -26602.001

$ws.Range("H72").Value2 = 5840
$ws.Range("I72").Value2 = 866.6667
$ws.Range("J72").Value2 = 8326.666999999999
$ws.Range("K72").Value2 = 7800.0003
$ws.Range("L72").Value2 = 74940.003
$ws.Range("M72").Value2 = -3744.0003
$ws.Range("N72").Value2 = -83052.003

$ws.Range("H100").Value2 = 2982.8572
$ws.Range("J100").Value2 = 2982.8572
$ws.Range("L100").Value2 = 8948.571599999999
$ws.Range("N100").Value2 = -10570.5716

$ws.Range("H113").Value2 = 2500645.5
$ws.Range("I113").Value2 = 555.75
$ws.Range("J113").Value2 = 6945249.5
$ws.Range("K113").Value2 = 1667.25
$ws.Range("L113").Value2 = 20835748.5
$ws.Range("M113").Value2 = 502.75
$ws.Range("N113").Value2 = -20840088.5

$ws = $wb.Worksheets("GSM")
$ws.Range("H107").Value2 = 729.5
$ws.Range("I107").Value2 = 549.5
$ws.Range("J107").Value2 = 999.5
$ws.Range("K107").Value2 = 549.5
$ws.Range("L107").Value2 = 999.5
$ws.Range("M107").Value2 = 1370.5
$ws.Range("N107").Value2 = -4839.5

$ws.Range("H122").Value2 = 11429.571
$ws.Range("I122").Value2 = 4969
$ws.Range("J122").Value2 = 16275
$ws.Range("K122").Value2 = 14907
$ws.Range("L122").Value2 = 48825
$ws.Range("M122").Value2 = -12457
$ws.Range("N122").Value2 = -53725

$ws = $wb.Worksheets("LTW")
$ws.Range("H122").Value2 = 8400
$ws.Range("I122").Value2 = 6500
$ws.Range("J122").Value2 = 8875
$ws.Range("K122").Value2 = 19500
$ws.Range("L122").Value2 = 26625
$ws.Range("M122").Value2 = -17050
$ws.Range("N122").Value2 = -31525

$ws.Range("H128").Value2 = 41941.332
$ws.Range("J128").Value2 = 41941.332
$ws.Range("L128").Value2 = 41941.332
$ws.Range("N128").Value2 = -51901.332

$ws.Range("H139").Value2 = 45078.332
$ws.Range("J139").Value2 = 45078.332
$ws.Range("L139").Value2 = 45078.332
$ws.Range("N139").Value2 = -55358.332

$ws = $wb.Worksheets("WVR")
$ws.Range("H122").Value2 = 10466.333
$ws.Range("I122").Value2 = 8900
$ws.Range("J122").Value2 = 11249.5
$ws.Range("K122").Value2 = 26700
$ws.Range("L122").Value2 = 33748.5
$ws.Range("M122").Value2 = -24250
$ws.Range("N122").Value2 = -38648.5

$ws.Range("H129").Value2 = 40725
$ws.Range("J129").Value2 = 40725
$ws.Range("L129").Value2 = 40725
$ws.Range("N129").Value2 = -50725
